# Updated cryptos list on Mon Apr 22 08:14:06 UTC 2024 with GitHub Actions
# Refresh Price (D) and Volume(1h) (E) columns with latest scrape; a few
# rows were also re-ranked (Coin name + Link swapped with a neighboring row).
# NumberFormat is forced to "@" (Text) immediately before writing any D-column
# value that looks like a plain number (e.g. "15.18"), otherwise Excel COM
# auto-coerces the string into a numeric cell -- then the cell's original
# style is restored so no stray formatting is introduced.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '65.993.44'
$ws.Range('E2').Value = '  +1.25%  '

$ws.Range('D3').Value = '3.201.05'
$ws.Range('E3').Value = '  +0.53%  '

$ws.Range('E4').Value = '  +0.01%  '

$origStyle = $ws.Range('D5').Style
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '602.88'
$ws.Range('D5').Style = $origStyle
$ws.Range('E5').Value = '  +3.76%  '

$origStyle = $ws.Range('D6').Style
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '153.52'
$ws.Range('D6').Style = $origStyle
$ws.Range('E6').Value = '  +1.31%  '

$ws.Range('E7').Value = '  -0.04%  '

$ws.Range('D8').Value = '3.197.99'
$ws.Range('E8').Value = '  +0.43%  '

$origStyle = $ws.Range('D9').Style
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.533'
$ws.Range('D9').Style = $origStyle
$ws.Range('E9').Value = '  -0.34%  '

$origStyle = $ws.Range('D10').Style
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.160'
$ws.Range('D10').Style = $origStyle
$ws.Range('E10').Value = '  -1.67%  '

$origStyle = $ws.Range('D11').Style
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.14'
$ws.Range('D11').Style = $origStyle
$ws.Range('E11').Value = '  -1.47%  '

$origStyle = $ws.Range('D12').Style
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.506'
$ws.Range('D12').Style = $origStyle
$ws.Range('E12').Value = '  -0.02%  '

$origStyle = $ws.Range('D13').Style
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000271'
$ws.Range('D13').Style = $origStyle
$ws.Range('E13').Value = '  -1.04%  '

$origStyle = $ws.Range('D14').Style
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '38.44'
$ws.Range('D14').Style = $origStyle
$ws.Range('E14').Value = '  +0.91%  '

$ws.Range('D15').Value = '3.728.25'
$ws.Range('E15').Value = '  +0.61%  '

$ws.Range('D16').Value = '66.130.53'
$ws.Range('E16').Value = '  +1.33%  '

$ws.Range('E17').Value = '  +3.30%  '

$ws.Range('D18').Value = '3.202.66'
$ws.Range('E18').Value = '  +0.29%  '

$ws.Range('E19').Value = '  -0.14%  '

$origStyle = $ws.Range('D20').Style
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '509.94'
$ws.Range('D20').Style = $origStyle
$ws.Range('E20').Value = '  -1.12%  '

$origStyle = $ws.Range('D21').Style
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '15.55'
$ws.Range('D21').Style = $origStyle
$ws.Range('E21').Value = '  +4.09%  '

$origStyle = $ws.Range('D22').Style
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.731'
$ws.Range('D22').Style = $origStyle
$ws.Range('E22').Value = '  -0.42%  '

$ws.Range('B23').Value = 'InternetComputer(DFINITY)'
$ws.Range('C23').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$origStyle = $ws.Range('D23').Style
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '15.18'
$ws.Range('D23').Style = $origStyle
$ws.Range('E23').Value = '  -0.52%  '

$ws.Range('B24').Value = 'Uniswap'
$ws.Range('C24').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$origStyle = $ws.Range('D24').Style
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '7.98'
$ws.Range('D24').Style = $origStyle
$ws.Range('E24').Value = '  +1.66%  '

$origStyle = $ws.Range('D25').Style
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '84.91'
$ws.Range('D25').Style = $origStyle
$ws.Range('E25').Value = '  -0.78%  '

$origStyle = $ws.Range('D26').Style
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.999'
$ws.Range('D26').Style = $origStyle
$ws.Range('E26').Value = '  -0.18%  '

$ws.Range('E27').Value = '  +2.02%  '

$origStyle = $ws.Range('D28').Style
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.15'
$ws.Range('D28').Style = $origStyle
$ws.Range('E28').Value = '  +1.20%  '

$ws.Range('E29').Value = '  +1.47%  '

$ws.Range('B30').Value = 'NEARProtocol'
$ws.Range('C30').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$origStyle = $ws.Range('D30').Style
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '6.86'
$ws.Range('D30').Style = $origStyle
$ws.Range('E30').Value = '  +8.06%  '

$ws.Range('B31').Value = 'Stacks'
$ws.Range('C31').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$origStyle = $ws.Range('D31').Style
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.86'
$ws.Range('D31').Style = $origStyle
$ws.Range('E31').Value = '  +2.36%  '

$origStyle = $ws.Range('D32').Style
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '28.07'
$ws.Range('D32').Style = $origStyle
$ws.Range('E32').Value = '  +0.05%  '

$origStyle = $ws.Range('D33').Style
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.21'
$ws.Range('D33').Style = $origStyle
$ws.Range('E33').Value = '  +0.61%  '

$origStyle = $ws.Range('D35').Style
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '6.59'
$ws.Range('D35').Style = $origStyle
$ws.Range('E35').Value = '  -0.74%  '

$ws.Range('E36').Value = '  -0.87%  '

$origStyle = $ws.Range('D37').Style
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0907'
$ws.Range('D37').Style = $origStyle
$ws.Range('E37').Value = '  -0.04%  '

$origStyle = $ws.Range('D38').Style
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '481.76'
$ws.Range('D38').Style = $origStyle
$ws.Range('E38').Value = '  +0.72%  '

$ws.Range('E39').Value = '  -0.66%  '

$origStyle = $ws.Range('D40').Style
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.95'
$ws.Range('D40').Style = $origStyle
$ws.Range('E40').Value = '  -6.55%  '

$origStyle = $ws.Range('D41').Style
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '8.81'
$ws.Range('D41').Style = $origStyle
$ws.Range('E41').Value = '  +1.18%  '

$origStyle = $ws.Range('D42').Style
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.297'
$ws.Range('D42').Style = $origStyle
$ws.Range('E42').Value = '  +3.32%  '

$ws.Range('E43').Value = '  -0.81%  '

$ws.Range('B44').Value = 'Maker'
$ws.Range('C44').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D44').Value = '2.944.31'
$ws.Range('E44').Value = '  -4.38%  '

$ws.Range('B45').Value = 'Fetch.AI'
$ws.Range('C45').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$origStyle = $ws.Range('D45').Style
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.46'
$ws.Range('D45').Style = $origStyle
$ws.Range('E45').Value = '  +1.21%  '

$ws.Range('D46').Value = '0.0₃0636'
$ws.Range('E46').Value = '  +3.62%  '

$origStyle = $ws.Range('D47').Style
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '28.66'
$ws.Range('D47').Style = $origStyle
$ws.Range('E47').Value = '  -2.10%  '

$ws.Range('E48').Value = '  +0.02%  '

$origStyle = $ws.Range('D49').Style
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.115'
$ws.Range('D49').Style = $origStyle
$ws.Range('E49').Value = '  -0.60%  '

$origStyle = $ws.Range('D50').Style
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.32'
$ws.Range('D50').Style = $origStyle
$ws.Range('E50').Value = '  +1.78%  '

$ws.Range('B51').Value = 'Monero'
$ws.Range('C51').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$origStyle = $ws.Range('D51').Style
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '119.80'
$ws.Range('D51').Style = $origStyle
$ws.Range('E51').Value = '  -1.01%  '
